# Applies the wk4.docx "sound changes + removed disused scenes" edit.
#
# Strategy: rewrite the body paragraph-by-paragraph using fresh Range
# objects. Lessons learned from probing this COM shim:
#   * A zero-width Range built fresh via $d.Range(n, n) at a position that
#     exactly coincides with a paragraph boundary binds to the *following*
#     paragraph, not the preceding one - so appends must be done with
#     paragraph.Range.InsertAfter(...) (using the paragraph's own Range
#     object), never via a separately constructed Range at Range.End.
#   * Deleting "Start..End+1" to swallow a paragraph mark is unsafe: an
#     already-empty paragraph's Range already spans its own mark (width 1
#     even though .Text reads ""), so "End+1" over-deletes into the next
#     paragraph. The robust way to delete paragraph N entirely (merging it
#     away) is Range(paragraph N start, paragraph N+1 start).Delete().
#   * Bookmarks.Add(name, range) anchored exactly at a point that further
#     InsertAfter calls subsequently touch will "ride along" with the
#     insertion and drift. So: insert all the text for a paragraph first,
#     *then* add the bookmark at the now-final, precise absolute offset.

$d = $word.ActiveDocument
$LF = [char]11   # Word "manual line break" (renders as <w:br/>)

function Set-ParaText($doc, $index, $text) {
    $p = $doc.Paragraphs.Item($index)
    $r = $doc.Range($p.Range.Start, $p.Range.End)
    $r.Text = $text
}

function Remove-Paragraph($doc, $index) {
    # Deletes paragraph $index (text + its own paragraph mark), merging it
    # away entirely - paragraph $index+1 becomes the new $index.
    $p = $doc.Paragraphs.Item($index)
    $pNext = $doc.Paragraphs.Item($index + 1)
    $delRange = $doc.Range($p.Range.Start, $pNext.Range.Start)
    $delRange.Delete()
}

# ---------------------------------------------------------------------
# 1) Remove the _GoBack bookmark from the "Date of Meeting" paragraph.
#    (It will be re-created later, anchored right after the new "Jam:"
#    run.)
# ---------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------
# 2) Paragraph 3: "Time of Meeting :" -> "Time of Meeting : 11:00"
# ---------------------------------------------------------------------
Set-ParaText $d 3 "Time of Meeting : 11:00"

# ---------------------------------------------------------------------
# 3) Paragraph 4: "Attendees:-" -> "Attendees:- Alex, Sean, Brittney, Josh"
# ---------------------------------------------------------------------
Set-ParaText $d 4 "Attendees:- Alex, Sean, Brittney, Josh"

# ---------------------------------------------------------------------
# 4) Paragraph 5: "Apologies from:-" -> "Apologies from:- n/a"
# ---------------------------------------------------------------------
Set-ParaText $d 5 "Apologies from:- n/a"

# Paragraph 6 ("Item One:-  Postmortem of previous week", bold) is unchanged.

# ---------------------------------------------------------------------
# 5) Paragraph 7: "What went well :  be specific " ->
#    "What went well :  visuals and mechanics almost ready for MVP pitch"
# ---------------------------------------------------------------------
Set-ParaText $d 7 "What went well :  visuals and mechanics almost ready for MVP pitch"

# ---------------------------------------------------------------------
# 6) Paragraph 8: "What went badly : be specific" ->
#    "What went badly : lack of progress with sprites / saved content"
# ---------------------------------------------------------------------
Set-ParaText $d 8 "What went badly : lack of progress with sprites / saved content"

# ---------------------------------------------------------------------
# 7) Paragraph 9: "Feedback Recieved : On any aspect ... etc." ->
#    "Feedback Recieved : audio relevant to moodboard, visuals tight and clean"
# ---------------------------------------------------------------------
Set-ParaText $d 9 "Feedback Recieved : audio relevant to moodboard, visuals tight and clean"

# Paragraph 10 ("Individual work completed:-") is unchanged.

# ---------------------------------------------------------------------
# 8) Paragraph 11: "Person 1<br>Person 2<br>Person 3<br>Person 4" ->
#    individual work update block, four lines separated by manual breaks.
# ---------------------------------------------------------------------
$text11 = "Alex: Mechanics and game engine work, audio iterated, implemented in engine, " + $LF + `
          "Sean: Art for world 1 complete, world 1 map in detail, bleed through assets, lighting, particle system, world 2 assets." + $LF + `
          "Brittney: Jellyfish concepts drawn" + $LF + `
          "Josh: research on playerprefs"
Set-ParaText $d 11 $text11

# ---------------------------------------------------------------------
# 9) Paragraph 12: "Item 2:-  Overall Aim ... X and Y" (partly bold) ->
#    "Item 2:-  MVP pitch will be complete by end of sprint." (no bold)
# ---------------------------------------------------------------------
Set-ParaText $d 12 "Item 2:-  MVP pitch will be complete by end of sprint."
$p12 = $d.Paragraphs.Item(12)
$r12 = $d.Range($p12.Range.Start, $p12.Range.End)
$r12.Bold = 0

# Paragraph 13 ("Tasks for the current week:-") is unchanged.

# ---------------------------------------------------------------------
# 10) Paragraph 14: "You need to be absolutely clear ... remotely delivered." ->
#     "Jam:" line (bookmark re-anchored right after it) followed by the
#     Alex+Josh / Sean / Brittney task lines.
# ---------------------------------------------------------------------
Set-ParaText $d 14 "Jam:"
$p14 = $d.Paragraphs.Item(14)
$jamBookmarkOffset = $p14.Range.Start + 4   # right after "Jam:"
$p14.Range.InsertAfter($LF + "Alex + Josh: Solve Playerprefs, add artifact triggers, " + $LF + `
          "Sean: Level 2 design" + $LF + `
          "Brittney: Complete character art")
$bmPoint = $d.Range($jamBookmarkOffset, $jamBookmarkOffset)
$d.Bookmarks.Add("_GoBack", $bmPoint)

# ---------------------------------------------------------------------
# 11) Paragraph 15: "Person 1 tasks<br>...<br>Person 4 tasks" ->
#     "Remote:" block.
# ---------------------------------------------------------------------
$text15 = "Remote:" + $LF + `
          "Alex: Sound effects" + $LF + `
          "Josh: build enemy behaviour" + $LF + `
          "Brittney: Complete Jellyfish Asset + Stationary enemy"
Set-ParaText $d 15 $text15

# ---------------------------------------------------------------------
# 12) Delete paragraph 16 "(These tasks to be uploaded and tracked on JIRA)"
#     and the following empty paragraph, entirely.
# ---------------------------------------------------------------------
Remove-Paragraph $d 16
Remove-Paragraph $d 16

# ---------------------------------------------------------------------
# 13) Paragraph 16 (now) "Item 3:-  Any Other Business." ->
#     "Item 3:-  " + line-break block about Level 1 / Sprites.
# ---------------------------------------------------------------------
Set-ParaText $d 16 "Item 3:-  "
$p16 = $d.Paragraphs.Item(16)
$p16.Range.InsertAfter($LF + "Level 1 is mostly complete, most mechanics are implemented, most level art is complete." + $LF + `
          "Sprites are critically needed as are sound effects ")

# ---------------------------------------------------------------------
# 14) Delete paragraph 17 (empty paragraph that used to separate AOB from
#     "Meeting Ended").
# ---------------------------------------------------------------------
Remove-Paragraph $d 17

# ---------------------------------------------------------------------
# 15) Paragraph 17 (now) "Meeting Ended :-" -> append " 12:00"
# ---------------------------------------------------------------------
$p17 = $d.Paragraphs.Item(17)
$p17.Range.InsertAfter(" 12:00")

# ---------------------------------------------------------------------
# 16) Paragraph 18 "Minute Taker:-" -> append " Alex"
# ---------------------------------------------------------------------
$p18 = $d.Paragraphs.Item(18)
$p18.Range.InsertAfter(" Alex")

# Paragraph 19 (trailing empty paragraph) is unchanged.

# ---------------------------------------------------------------------
# Dump the final structure for verification.
# ---------------------------------------------------------------------
Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
$i = 1
foreach ($p in $d.Paragraphs) {
    Write-Output "$i : [$($p.Range.Text)]"
    $i = $i + 1
}
